# Automatic update of files.
# This script re-applies the "Id" (and associated record) shuffle that
# occurred across rows 9, 11, 12, 14, 15, 18, 19, 20, 21, 22 of the sheet.
# Each of those rows effectively receives the full content of another one
# of those rows (a permutation), while every other cell/row is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source row number (content to copy from)
$mapping = @{
    9  = 21
    11 = 18
    12 = 15
    14 = 9
    15 = 11
    18 = 19
    19 = 20
    20 = 12
    21 = 22
    22 = 14
}

# Snapshot the "before" values of every involved row/column, since several
# rows both give and receive data and we must not read already-overwritten
# cells while computing later assignments.
$cols = @("A","B","D","E","F","G","H","Q","R","AC","AJ","AK","AO")
$rows = @(9,11,12,14,15,18,19,20,21,22)

$before = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $addr = "$c$r"
        $rowVals[$c] = $ws.Range($addr).Value()
    }
    $before[$r] = $rowVals
}

# Apply the new values (taken from the snapshot of the mapped source row)
foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $addr = "$c$newRow"
        $val = $srcVals[$c]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
